# Applies the "Added README tab + added meteo to functionalities" commit.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Source table: flip the "* Meteo" row (row 28) to supported, and
#    fill in the Supported-since / module / class columns. The three
#    "<topic> mkdocs table" sheets spill from a FILTER array formula
#    that reads this table, so they recompute automatically.
# ------------------------------------------------------------------
$wsSource = $wb.Worksheets.Item("Source table")
$wsSource.Range("C28").Value = "X"
$wsSource.Range("E28").Value = "X"
$wsSource.Range("G28").Value = "0.5.0"
$wsSource.Range("H28").Value = "hydrolib.core.dflowfm.ext.models"
$wsSource.Range("I28").Value = "Meteo"

# ------------------------------------------------------------------
# 2. Add the new "README" tab at the end of the workbook (after
#    "Topics") describing usage / inner workings of the sheet.
# ------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsReadme = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsReadme.Name = "README"

$wsReadme.Range("A1").Value = "Usage"
$wsReadme.Range("A1").Font.Bold = $true
$wsReadme.Range("A2").Value = "* Adding/changing entries in the functionalities list should be done only in the Source table tab sheet."
$wsReadme.Range("A3").Value = "* Only edit the visible (non-Hidden) columns. The two hidden columns automatically produce Markdown code for the table icons."

$wsReadme.Range("A5").Value = "Inner workings"
$wsReadme.Range("A5").Font.Bold = $true
$wsReadme.Range("A6").Value = "* The three first tab sheets automatically fill up based on the Source table and the kernel topic name (using FILTER and SELECTCOLS)"
$wsReadme.Range("A7").Value = "* As the Source table grows longer, make sure that on the first three tabs the select formula in A3 still includes all row numbers for which Source table has data values."

$wsReadme.Range("A8").Select() | Out-Null

# ------------------------------------------------------------------
# 3. Restore / update the per-sheet selections (and DIMR's
#    "show formulas" view) left behind by the author's editing
#    session, finishing on "Source table" so it stays the active tab.
# ------------------------------------------------------------------
$wsDIMR = $wb.Worksheets.Item("DIMR mkdocs table")
$wsDIMR.Activate() | Out-Null
$excel.ActiveWindow.DisplayFormulas = $true
$wsDIMR.Range("A3").Select() | Out-Null

$wsRR = $wb.Worksheets.Item("RR mkdocs table")
$wsRR.Activate() | Out-Null
$wsRR.Range("A2").Select() | Out-Null

$wsFM = $wb.Worksheets.Item("FM mkdocs table")
$wsFM.Activate() | Out-Null
$wsFM.Range("A3").Select() | Out-Null

$wsSource.Activate() | Out-Null
$wsSource.Range("C37").Select() | Out-Null
